# #441 Added style column
#
# Inserts a new "Style" column into the Products sheet (between "Brand Name"
# and "Seller"), populates the header + sample values, and leaves the
# Products sheet as the active/selected tab (it previously was the Images
# sheet that was left selected).

$wb = $excel.ActiveWorkbook

$products = $wb.Worksheets.Item("Products")

# Insert a new column at J (10th column), shifting Seller..Quantity one
# column to the right (J:N -> K:O).
$products.Columns.Item(10).Insert()

# New "Style" column header + data.
$products.Range("J1").Value = "Style"
$products.Range("J2").Value = 1
$products.Range("J3").Value = 1

# Make the Products sheet the active tab/selection (previously Images was
# the selected tab).
$products.Activate()
$products.Range("G15").Select()
